$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.601837382907111, 1.6659946294738037, -1.024585524996625, 2.5351015952785687)"
$ws.Range("C2").Value = "MIE(8.023361018571332, 4.008816126641602, -9.676324522946352, 14.76630828651298)"
$ws.Range("D2").Value = "JSU(-0.8201033655273211, 1.0382123968132873, 0.5147207733496624, 2.8291306301441965)"
$ws.Range("E2").Value = "JSU(-0.7905148869504085, 1.078562365629914, 4.771277283980243, 4.228684894806838)"
